$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("baseline")
$ws.Activate()

# Add new row of data for "roberta + bilstm"
$ws.Range("B12").Value = "roberta + bilstm"
$ws.Range("C12").Value = 1.5179
$ws.Range("D12").Value = 0.6692
$ws.Range("E12").Value = 0.671
$ws.Range("F12").Value = 0.6692
$ws.Range("G12").Value = 0.6559
$ws.Range("H12").Value = 0.3739
$ws.Range("I12").Value = 0.3792
$ws.Range("J12").Value = 0.5673

# Match style (number format) of the rest of the data rows
$ws.Range("C12:J12").NumberFormat = "0.0000"

# Update column B width to fit the new, longer text (best-fit to contents,
# matching the "bestFit" auto-sized behavior applied by Excel for this column)
$ws.Columns.Item(2).AutoFit() | Out-Null

# Update the active cell selection as shown in the diff
$ws.Range("I14").Select()
